$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New English translations in column D for rows 29-41 (dialogue lines)
$ws.Range("D29").Value = "\n<リリー>Oh Shina, you're back.`nDid you lock them up?`n"
$ws.Range("D30").Value = "\n<シィナ>Huh, me?`nBeats me."
$ws.Range("D31").Value = "\n<リリー>Whoa whoa whoa!`nYou need to take this seriously!`nWhat if they got out!"
$ws.Range("D32").Value = "\n<リリー>Sigh....`nWhy are you always so... so..."
$ws.Range("D33").Value = "\n<シィナ>Whaa?`nI dunno-nya.`nI don't have the key. It's not my fault!"
$ws.Range("D34").Value = "\n<ライム>I don't have a key either."
$ws.Range("D35").Value = "\n<リリー>Hm?"
$ws.Range("D36").Value = "\n<リリー>...!!"
$ws.Range("D37").Value = "\n<リリー>I have the key.`n...`nDid you at least shut the door?"
$ws.Range("D38").Value = "\n<シィナ>I dunno.`nWhy don't you do it yourself-nya."
$ws.Range("D39").Value = "\n<ライム>Well for now, let's go back and check,`nLily?"
$ws.Range("D40").Value = "\n<シィナ>Now who should be the one apologizing?~"
$ws.Range("D41").Value = "\n<リリー>Eei♥"

# Column D duplicates column A content for rows 42-52 (non-dialogue / already-shared strings)
$ws.Range("D42").Value = "最初のイベント"
$ws.Range("D43").Value = "\n<ライム>やっぱりね！"
$ws.Range("D44").Value = "\n<ライム>逃げる時間そんなにないよなー、って思ってたの。`nどこかに隠れてたのかな？`n待ってたら来ると思ったよー！"
$ws.Range("D45").Value = "\n<ライム>にへへへへー♥`nつーかまーえた♥"
$ws.Range("D46").Value = "\n<\n[3]>ぬるぬるして動きにくいでしょー。`n早く逃げないと白いの出させちゃうぞー？"
$ws.Range("D47").Value = "\C[3]※捕まるとタイミングバーが表示されます。`n\C[0]タイミングよく黄か赤で止めてください。`n赤で止めると被ダメージが半減します。"
$ws.Range("D48").Value = "\n<\n[3]>むにゅー・・・♥`n柔らかくて溶けちゃいそうでしょー♥`n気持ちいい気持ちいいー♥"
$ws.Range("D49").Value = "\n<\n[3]>あれー？もう出ちゃうのー？`nおっぱい我慢できなかったー？あはは♥`nじゃあ一回だけ、出しちゃおっかー♪"
$ws.Range("D50").Value = "\n<\n[3]>あっあっ♥おっぱいの間でぴくぴくしてるー♥`nおちんちん喜んでくれたみたい♥`nうれしー♥"
$ws.Range("D51").Value = "\C[1]SAN値が1下がった・・・（現在SAN値\v[270]）"
$ws.Range("D52").Value = "\n<ライム>もうー。ひょっとして全然抵抗する気ないのー？`nゲームはまだ始まったばかりなのに。`n次は本気で搾っちゃうからねー？"

# Reset auto row-height override introduced by multi-line cell entry so rows match default height
$ws.Rows.Item(29).AutoFit()
$ws.Rows.Item(30).AutoFit()
$ws.Rows.Item(31).AutoFit()
$ws.Rows.Item(32).AutoFit()
$ws.Rows.Item(33).AutoFit()
$ws.Rows.Item(34).AutoFit()
$ws.Rows.Item(35).AutoFit()
$ws.Rows.Item(36).AutoFit()
$ws.Rows.Item(37).AutoFit()
$ws.Rows.Item(38).AutoFit()
$ws.Rows.Item(39).AutoFit()
$ws.Rows.Item(40).AutoFit()
$ws.Rows.Item(41).AutoFit()
$ws.Rows.Item(42).AutoFit()
$ws.Rows.Item(43).AutoFit()
$ws.Rows.Item(44).AutoFit()
$ws.Rows.Item(45).AutoFit()
$ws.Rows.Item(46).AutoFit()
$ws.Rows.Item(47).AutoFit()
$ws.Rows.Item(48).AutoFit()
$ws.Rows.Item(49).AutoFit()
$ws.Rows.Item(50).AutoFit()
$ws.Rows.Item(51).AutoFit()
$ws.Rows.Item(52).AutoFit()
